$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in shared string used by A1:A4 ("Massachusettes" -> "Massachusetts")
$ws.Range("A1:A4").Value = "Massachusetts {{ ma }}"

# Add new data row 10 with numeric values
$ws.Range("A10").Value = 2019
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 16

# Update the active selection to match the saved view state
$ws.Range("B15").Select()
